$d = $word.ActiveDocument

# Apply each replacement in document order (top to bottom) so that
# newly-inserted text from an earlier replacement is never re-matched
# by a later rule (e.g. "24÷8=" -> "66÷7=" happens after
# "66÷7=" -> "99÷4=" has already run).

$d.Content.Find.Execute("27÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "45÷8=", 2)
$d.Content.Find.Execute("58÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "33÷9=", 2)
$d.Content.Find.Execute("22÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "92÷4=", 2)
$d.Content.Find.Execute("50÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "60÷2=", 2)
$d.Content.Find.Execute("16÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "14÷8=", 2)
$d.Content.Find.Execute("74÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "50÷2=", 2)
$d.Content.Find.Execute("24÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "88÷9=", 2)
$d.Content.Find.Execute("66÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "99÷4=", 2)
$d.Content.Find.Execute("19÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "34÷8=", 2)
$d.Content.Find.Execute("39÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "76÷4=", 2)
$d.Content.Find.Execute("24÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "84÷9=", 2)
$d.Content.Find.Execute("55÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "40÷7=", 2)
$d.Content.Find.Execute("77÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "81÷2=", 2)
$d.Content.Find.Execute("48÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "26÷9=", 2)
$d.Content.Find.Execute("99÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "11÷5=", 2)
$d.Content.Find.Execute("24÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "66÷7=", 2)
$d.Content.Find.Execute("26÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "68÷2=", 2)
$d.Content.Find.Execute("96÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "60÷8=", 2)
$d.Content.Find.Execute("62÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "76÷9=", 2)
$d.Content.Find.Execute("95÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "90÷2=", 2)
$d.Content.Find.Execute("29÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "17÷5=", 2)
$d.Content.Find.Execute("68÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "92÷5=", 2)
$d.Content.Find.Execute("85÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "13÷9=", 2)
$d.Content.Find.Execute("11÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "93÷6=", 2)
$d.Content.Find.Execute("79÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "32÷3=", 2)
